$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates: force text to preserve formatting, avoid numeric coercion ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.658.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.597.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.514"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0619"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.248"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.822.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.672.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.635.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "208.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.114"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0515"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.282.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.618"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.49"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.734.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.102"
$ws.Range("D48").Style = "Normal"

# --- Column B/C/E updates ---
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("E13").Value = "  +4.58%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("E21").Value = "  +3.68%  "
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("E23").Value = "  +2.84%  "
$ws.Range("E24").Value = "  +0.66%  "
$ws.Range("E25").Value = "  -1.33%  "
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("E27").Value = "  -1.11%  "
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("E34").Value = "  -1.15%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E35").Value = "  -7.37%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("E38").Value = "  -0.82%  "
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("E40").Value = "  +17.26%  "
$ws.Range("E41").Value = "  +1.94%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("E47").Value = "  -3.52%  "
$ws.Range("E48").Value = "  +1.27%  "
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("E51").Value = "  -1.99%  "
